$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto table: Price (D) and Volume(1h) (E) columns for rows 2-51,
# matching the "Updated cryptos list ... with GitHub Actions" feed refresh.
#
# The source feed always writes these two columns as plain text (e.g. the
# Price column holds locale-formatted strings like "28.244.08" that are not
# valid numbers, and Volume(1h) keeps its "  +0.80%  " padding) so the whole
# D2:E51 block is pre-formatted as Text before writing; this stops Excel from
# auto-converting purely-numeric-looking values (e.g. "1.004", "0.9999") into
# real numbers. The NumberFormat is reset back to Normal afterwards so no
# visible formatting change is left behind on the cells.
$dataRng = $ws.Range("D2:E51")
$dataRng.NumberFormat = "@"

$ws.Range("D2").Value = "28.244.08"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.803.82"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "337.12"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.4664"
$ws.Range("E7").Value = "  +21.86%  "
$ws.Range("D8").Value = "0.3798"
$ws.Range("E8").Value = "  +10.44%  "
$ws.Range("D9").Value = "45.12"
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "0.07629"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("D11").Value = "1.149"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "22.40"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "6.331"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "7.466"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "1.803.03"
$ws.Range("E16").Value = "  +2.96%  "
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "0.06736"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("E19").Value = "  +3.15%  "
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("D22").Value = "6.405"
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("D23").Value = "28.245.58"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "11.87"
$ws.Range("E24").Value = "  +1.48%  "
$ws.Range("D25").Value = "2.409"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +4.40%  "
$ws.Range("D27").Value = "153.80"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "2.371"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").Value = "2.010.34"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").Value = "133.67"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "4.037"
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "0.09584"
$ws.Range("E33").Value = "  +8.78%  "
$ws.Range("D34").Value = "5.856"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "0.2228"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").Value = "0.06365"
$ws.Range("E36").Value = "  +3.22%  "
$ws.Range("D37").Value = "12.10"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "0.02350"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").Value = "5.260"
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").Value = "0.6633"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("D43").Value = "8.214"
$ws.Range("E43").Value = "  +3.06%  "
$ws.Range("D44").Value = "14.24"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.6119"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "3.835"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "130.18"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("E51").Value = "  +0.72%  "

$dataRng.Style = "Normal"
